# Update "想去人数" (want-to-go counts) on the two sheets that list
# full event data: "展览" (sheet 1) and "全部类型" (sheet 4).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F7").Value = 347
$ws1.Range("F8").Value = 4937
$ws1.Range("F10").Value = 5224
$ws1.Range("F11").Value = 600
$ws1.Range("F12").Value = 1311

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F8").Value = 347
$ws4.Range("F9").Value = 4937
$ws4.Range("F11").Value = 5224
$ws4.Range("F12").Value = 600
$ws4.Range("F13").Value = 1311
